$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ECs" sending-cluster row (old row 2); FAPs/MuSCs rows shift up.
$ws.Rows("2:2").Delete()

# Recalculated derived-specificity values (now that ECs is excluded from the set).
$ws.Range("I2").Value = 0.1949338371837906
$ws.Range("J2").Value = 0.1949338371837907
$ws.Range("S2").Value = 0.1949338371837906
$ws.Range("T2").Value = 0.1949338371837907

$ws.Range("I3").Value = 0.8050661628162092
$ws.Range("J3").Value = 0.8050661628162094
$ws.Range("S3").Value = 0.8050661628162092
$ws.Range("T3").Value = 0.8050661628162094
